$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.715.74'
$ws.Range('E2').Value = '  +0.49%  '
$ws.Range('D3').Value = '1.599.47'
$ws.Range('E3').Value = '  +1.01%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '211.24'
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('E6').Value = '  +1.33%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.0618'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('E9').Value = '  -1.21%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.62'
$ws.Range('E10').Value = '  +0.61%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0844'
$ws.Range('E11').Value = '  +1.30%  '
$ws.Range('B12').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C12').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D12').Value = '1.823.63'
$ws.Range('E12').Value = '  +0.99%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.632.05'
$ws.Range('E13').Value = '  +3.41%  '
$ws.Range('E14').Value = '  -0.49%  '
$ws.Range('E15').Value = '  -1.14%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.92'
$ws.Range('E16').Value = '  +0.56%  '
$ws.Range('D17').Value = '26.690.25'
$ws.Range('E17').Value = '  +0.32%  '
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '208.53'
$ws.Range('E19').Value = '  +0.48%  '
$ws.Range('E20').Value = '  -0.03%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.80'
$ws.Range('E21').Value = '  +0.55%  '
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.35'
$ws.Range('E23').Value = '  -1.01%  '
$ws.Range('E24').Value = '  -0.22%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '145.64'
$ws.Range('E25').Value = '  -0.56%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.24'
$ws.Range('E27').Value = '  -2.24%  '
$ws.Range('E28').Value = '  +1.83%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0507'
$ws.Range('E30').Value = '  +0.92%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.16'
$ws.Range('E31').Value = '  +0.41%  '
$ws.Range('E32').Value = '  -0.66%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.665'
$ws.Range('E33').Value = '  -2.86%  '
$ws.Range('E34').Value = '  +0.51%  '
$ws.Range('D35').Value = '1.288.37'
$ws.Range('E35').Value = '  -2.44%  '
$ws.Range('E36').Value = '  -1.10%  '
$ws.Range('E37').Value = '  -0.66%  '
$ws.Range('E38').Value = '  -0.37%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.844'
$ws.Range('E39').Value = '  +2.58%  '
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.43'
$ws.Range('E41').Value = '  +1.74%  '
$ws.Range('E42').Value = '  +1.03%  '
$ws.Range('E43').Value = '  +0.43%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '63.66'
$ws.Range('E44').Value = '  +0.37%  '
$ws.Range('D45').Value = '1.736.18'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.905'
$ws.Range('E46').Value = '  +8.75%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '90.08'
$ws.Range('E47').Value = '  +0.60%  '
$ws.Range('E48').Value = '  -0.58%  '
$ws.Range('E49').Value = '  +2.50%  '
$ws.Range('E50').Value = '  -0.28%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '7.49'
$ws.Range('E51').Value = '  +0.00%  '
